# Apply the update described by the diff:
#  - Fix a typo in the existing "ödeme" scenario row (row 44, column B):
#      "tutarlları" -> "tutarları"
#  - Append two new scenario rows about "ihale" (tender) to the bottom of
#    the table (rows 45 and 46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in existing row 44 (Senaryo / column B) ---
$ws.Range("B44").Value = "Müşterinin ödediği tutarları görmek istiyorum."

# --- New row 45: ihale / listeleme kanalı ---
$ws.Range("A45").Value = "ihale"
$ws.Range("B45").Value = "Aracın hangi ihaleden satışa çıktığını görmek istiyorum"
$ws.Range("C45").Value = "Fırsat ekranında üstte Listelendiği kanal alanı var orada görebilirsin."
$ws.Range("D45").Value = "Salesforce da plakayı yazıp ara. Ödeme - Evrak - Teslimat- Kazanıldı aşamasında olan kaydı seç. Kalın punto fırsat başlığının hemen altında Listedenliği kanal olarak göreceksin."
$ws.Range("E45").Value = "Product Manager"
$ws.Range("F45").Value = "listeleme kanalı.JPG"

# --- New row 46: ihale / ihale bitiş tarihi ---
$ws.Range("A46").Value = "ihale"
$ws.Range("B46").Value = "Aracın ihale bitiş tarihini görmek istiyorum"
$ws.Range("C46").Value = "Fırsat ekranında en altta solda oluşturma tarihi var. Orada görebilirsin."
$ws.Range("D46").Value = "Salesforce da plakayı yazıp ara. Ödeme - Evrak - Teslimat- Kazanıldı aşamasında olan kaydı seç. En alta in. Bu ihale bittikten sonra gelen fırsatın oluşturulduğu tarihtir."
$ws.Range("E46").Value = "Product Manager"
$ws.Range("F46").Value = "ihale tarihi.JPG"

# Update the active selection to mirror the author's post-edit state
# (selection below the newly-added data, matching the diff).
$ws.Range("F47").Select()
